$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.007.51"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "3.268.85"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'575.26"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'179.56"
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  +3.14%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("D10").Value = "'6.69"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").Value = "'0.399"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").Value = "3.845.82"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("E13").Value = "  -3.73%  "
$ws.Range("D14").Value = "66.087.73"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").Value = "'26.33"
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").Value = "3.269.56"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").Value = "'429.49"
$ws.Range("E18").Value = "  -3.07%  "
$ws.Range("D19").Value = "'5.51"
$ws.Range("E19").Value = "  -2.61%  "
$ws.Range("D20").Value = "'13.05"
$ws.Range("E20").Value = "  -3.57%  "
$ws.Range("D21").Value = "'7.35"
$ws.Range("E21").Value = "  -4.48%  "
$ws.Range("D22").Value = "'71.65"
$ws.Range("E22").Value = "  -3.26%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "3.418.17"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").Value = "'0.502"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("D28").Value = "'8.79"
$ws.Range("E28").Value = "  -2.55%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("D31").Value = "'22.15"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "'5.13"
$ws.Range("E33").Value = "  -3.62%  "
$ws.Range("D34").Value = "'6.54"
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("E35").Value = "  -4.32%  "
$ws.Range("D36").Value = "'157.78"
$ws.Range("E36").Value = "  -2.82%  "
$ws.Range("D37").Value = "'1.41"
$ws.Range("E37").Value = "  -5.83%  "
$ws.Range("D38").Value = "'26.32"
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("E39").Value = "  -3.31%  "
$ws.Range("D40").Value = "2.752.76"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").Value = "'0.773"
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("D42").Value = "'4.27"
$ws.Range("E42").Value = "  -4.29%  "
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  -3.10%  "
$ws.Range("D45").Value = "'0.0653"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").Value = "'2.29"
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("D47").Value = "'319.19"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").Value = "'23.04"
$ws.Range("E48").Value = "  -6.22%  "
$ws.Range("D49").Value = "'0.0264"
$ws.Range("E49").Value = "  -2.65%  "
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("E51").Value = "  +0.05%  "
